$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "prodXTCt"
$ws.Range("B5").Value = "prodJWne"
$ws.Range("B2").Value = "prodQPJL"
